# Update inverter capacity results with freshly recomputed values from server.
$wb = $excel.ActiveWorkbook

# --- Sheet "2025" ---
$ws = $wb.Worksheets.Item("2025")
$ws.Range("O2").Value = 3.41272695030143

# --- Sheet "2030" ---
$ws = $wb.Worksheets.Item("2030")
$ws.Range("B2").Value = 0.06592367743582367
$ws.Range("I2").Value = 0.7356952478631865
$ws.Range("L2").Value = 0.2195099354701471
$ws.Range("M2").Value = 0.08791305939480454
$ws.Range("N2").Value = 8.998982253460907
$ws.Range("O2").Value = 6.747238139072843

# --- Sheet "2035" ---
$ws = $wb.Worksheets.Item("2035")
$ws.Range("A2").Value = 0.1260055233262948
$ws.Range("B2").Value = 0.03282537593862538
$ws.Range("E2").Value = 0.1576050149385859
$ws.Range("I2").Value = 0.4124253487523346
$ws.Range("M2").Value = 0.02894162393852881
$ws.Range("N2").Value = 7.791100800027102
$ws.Range("O2").Value = 0.7270200098146926

# --- Sheet "2040" ---
$ws = $wb.Worksheets.Item("2040")
$ws.Range("N2").Value = 0.4711113629604995

# --- Sheet "2045" ---
$ws = $wb.Worksheets.Item("2045")
$ws.Range("A2").Value = 0.1538638740281748
$ws.Range("N2").Value = 2.534564915429591
$ws.Range("O2").Value = 5.522246454485551
